$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1887316  # H17: was 1961252.2
$ws.Cells.Item(17, 10).Value = 1923597.8  # J17: was 2000464
$ws.Cells.Item(17, 12).Value = 5770793.4  # L17: was 6001392
$ws.Cells.Item(17, 14).Value = -5771129.4  # N17: was -6001728
$ws.Cells.Item(40, 8).Value = 1413.3  # H40: was 2066.6
$ws.Cells.Item(40, 9).Value = 1257.1428  # I40: was 2000
$ws.Cells.Item(40, 10).Value = 1777.6666  # J40: was 2166.5
$ws.Cells.Item(40, 11).Value = 1257.1428  # K40: was 2000
$ws.Cells.Item(40, 12).Value = 1777.6666  # L40: was 2166.5
$ws.Cells.Item(40, 13).Value = -1082.1428  # M40: was -1825
$ws.Cells.Item(40, 14).Value = -2127.6666  # N40: was -2516.5
$ws.Cells.Item(41, 8).Value = 324.82352  # H41: was 300.6842
$ws.Cells.Item(41, 9).Value = 102.5  # I41: was 98.666664
$ws.Cells.Item(41, 10).Value = 354.46667  # J41: was 338.5625
$ws.Cells.Item(41, 11).Value = 102.5  # K41: was 98.666664
$ws.Cells.Item(41, 12).Value = 354.46667  # L41: was 338.5625
$ws.Cells.Item(41, 13).Value = 337.5  # M41: was 341.333336
$ws.Cells.Item(41, 14).Value = -1234.46667  # N41: was -1218.5625
$ws.Cells.Item(70, 8).Value = 1057.2858  # H70: was 1146.3334
$ws.Cells.Item(70, 9).Value = 1045.6364  # I70: was 1102.125
$ws.Cells.Item(70, 10).Value = 1100  # J70: was 1500
$ws.Cells.Item(70, 11).Value = 3136.9092  # K70: was 3306.375
$ws.Cells.Item(70, 12).Value = 3300  # L70: was 4500
$ws.Cells.Item(70, 13).Value = -2866.9092  # M70: was -3036.375
$ws.Cells.Item(70, 14).Value = -3840  # N70: was -5040
$ws.Cells.Item(73, 8).Value = 1057.2858  # H73: was 1146.3334
$ws.Cells.Item(73, 9).Value = 1045.6364  # I73: was 1102.125
$ws.Cells.Item(73, 10).Value = 1100  # J73: was 1500
$ws.Cells.Item(73, 11).Value = 3136.9092  # K73: was 3306.375
$ws.Cells.Item(73, 12).Value = 3300  # L73: was 4500
$ws.Cells.Item(73, 13).Value = -2200.9092  # M73: was -2370.375
$ws.Cells.Item(73, 14).Value = -5172  # N73: was -6372
$ws.Cells.Item(100, 8).Value = 3165  # H100: was 2915
$ws.Cells.Item(100, 9).Value = 2372.5  # I100: was 2498
$ws.Cells.Item(100, 10).Value = 4750  # J100: was 5000
$ws.Cells.Item(100, 11).Value = 2372.5  # K100: was 2498
$ws.Cells.Item(100, 12).Value = 4750  # L100: was 5000
$ws.Cells.Item(100, 13).Value = -1831.5  # M100: was -1957
$ws.Cells.Item(100, 14).Value = -5832  # N100: was -6082
$ws.Cells.Item(111, 8).Value = 6646.857  # H111: was 5048.6
$ws.Cells.Item(111, 9).Value = 2514  # I111: was 1797.2
$ws.Cells.Item(111, 11).Value = 7542  # K111: was 5391.6
$ws.Cells.Item(111, 13).Value = -4475  # M111: was -2324.6
$ws.Cells.Item(113, 8).Value = 19233700  # H113: was 20411202
$ws.Cells.Item(113, 9).Value = 28573310  # I113: was 31251950
$ws.Cells.Item(113, 10).Value = 5087.8823  # J113: was 5088.1763
$ws.Cells.Item(113, 11).Value = 28573310  # K113: was 31251950
$ws.Cells.Item(113, 12).Value = 5087.8823  # L113: was 5088.1763
$ws.Cells.Item(113, 13).Value = -28570056  # M113: was -31248696
$ws.Cells.Item(113, 14).Value = -11595.8823  # N113: was -11596.1763
$ws.Cells.Item(129, 8).Value = 147896.19  # H129: was 162126.05
$ws.Cells.Item(129, 10).Value = 154682.33  # J129: was 170325.84
$ws.Cells.Item(129, 12).Value = 464046.99  # L129: was 510977.52
$ws.Cells.Item(129, 14).Value = -474046.99  # N129: was -520977.52
$ws.Cells.Item(138, 8).Value = 1528.697  # H138: was 1401.8387
$ws.Cells.Item(138, 10).Value = 3244.5  # J138: was 3194.4
$ws.Cells.Item(138, 12).Value = 9733.5  # L138: was 9583.200000000001
$ws.Cells.Item(138, 14).Value = -20013.5  # N138: was -19863.2

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3348.3684  # H45: was 3298.3845
$ws.Cells.Item(45, 10).Value = 3246.4167  # J45: was 3172.52
$ws.Cells.Item(45, 12).Value = 3246.4167  # L45: was 3172.52
$ws.Cells.Item(45, 14).Value = -4000.4167  # N45: was -3926.52
$ws.Cells.Item(61, 8).Value = 3222.818  # H61: was 3222.8635
$ws.Cells.Item(61, 9).Value = 2218.875  # I61: was 2218.9375
$ws.Cells.Item(61, 11).Value = 2218.875  # K61: was 2218.9375
$ws.Cells.Item(61, 13).Value = -2006.875  # M61: was -2006.9375
$ws.Cells.Item(97, 8).Value = 1920  # H97: was 1807.8572
$ws.Cells.Item(97, 9).Value = 1790  # I97: was 1669.1666
$ws.Cells.Item(97, 10).Value = 2635  # J97: was 2640
$ws.Cells.Item(97, 11).Value = 1790  # K97: was 1669.1666
$ws.Cells.Item(97, 12).Value = 2635  # L97: was 2640
$ws.Cells.Item(97, 13).Value = -1294  # M97: was -1173.1666
$ws.Cells.Item(97, 14).Value = -3627  # N97: was -3632
$ws.Cells.Item(110, 8).Value = 1573.4762  # H110: was 1802.9445
$ws.Cells.Item(110, 9).Value = 2060.1  # I110: was 2526.375
$ws.Cells.Item(110, 10).Value = 1131.091  # J110: was 1224.2
$ws.Cells.Item(110, 11).Value = 2060.1  # K110: was 2526.375
$ws.Cells.Item(110, 12).Value = 1131.091  # L110: was 1224.2
$ws.Cells.Item(110, 13).Value = -15.09999999999991  # M110: was -481.375
$ws.Cells.Item(110, 14).Value = -5221.091  # N110: was -5314.2
$ws.Cells.Item(136, 8).Value = 3222.818  # H136: was 3222.8635
$ws.Cells.Item(136, 9).Value = 2218.875  # I136: was 2218.9375
$ws.Cells.Item(136, 11).Value = 6656.625  # K136: was 6656.8125
$ws.Cells.Item(136, 13).Value = -4106.625  # M136: was -4106.8125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2277.5454  # H86: was 2304.762
$ws.Cells.Item(86, 9).Value = 1877.3334  # I86: was 1887.4117
$ws.Cells.Item(86, 11).Value = 1877.3334  # K86: was 1887.4117
$ws.Cells.Item(86, 13).Value = -754.3334  # M86: was -764.4117000000001
$ws.Cells.Item(89, 8).Value = 2277.5454  # H89: was 2304.762
$ws.Cells.Item(89, 9).Value = 1877.3334  # I89: was 1887.4117
$ws.Cells.Item(89, 11).Value = 9386.666999999999  # K89: was 9437.058500000001
$ws.Cells.Item(89, 13).Value = -3770.666999999999  # M89: was -3821.058500000001
$ws.Cells.Item(94, 8).Value = 2923.0293  # H94: was 3469.75
$ws.Cells.Item(94, 9).Value = 1462.5  # I94: was 1638.579
$ws.Cells.Item(94, 10).Value = 5600.6665  # J94: was 7335.5557
$ws.Cells.Item(94, 11).Value = 1462.5  # K94: was 1638.579
$ws.Cells.Item(94, 12).Value = 5600.6665  # L94: was 7335.5557
$ws.Cells.Item(94, 13).Value = -1011.5  # M94: was -1187.579
$ws.Cells.Item(94, 14).Value = -6502.6665  # N94: was -8237.555700000001
$ws.Cells.Item(99, 8).Value = 2212.7144  # H99: was 2225.5715
$ws.Cells.Item(99, 9).Value = 2061.8  # I99: was 1766.5
$ws.Cells.Item(99, 10).Value = 2590  # J99: was 4980
$ws.Cells.Item(99, 11).Value = 2061.8  # K99: was 1766.5
$ws.Cells.Item(99, 12).Value = 2590  # L99: was 4980
$ws.Cells.Item(99, 13).Value = -563.8000000000002  # M99: was -268.5
$ws.Cells.Item(99, 14).Value = -5586  # N99: was -7976

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1163.1666  # H22: was 623.75
$ws.Cells.Item(22, 9).Value = 899.8  # I22: was 450
$ws.Cells.Item(22, 10).Value = 2480  # J22: was 913.3333
$ws.Cells.Item(22, 11).Value = 899.8  # K22: was 450
$ws.Cells.Item(22, 12).Value = 2480  # L22: was 913.3333
$ws.Cells.Item(22, 13).Value = -549.8  # M22: was -100
$ws.Cells.Item(22, 14).Value = -3180  # N22: was -1613.3333
$ws.Cells.Item(99, 8).Value = 5513.1875  # H99: was 4907302
$ws.Cells.Item(99, 9).Value = 4246.4546  # I99: was 9263126
$ws.Cells.Item(99, 10).Value = 8300  # J99: was 7000
$ws.Cells.Item(99, 11).Value = 4246.4546  # K99: was 9263126
$ws.Cells.Item(99, 12).Value = 8300  # L99: was 7000
$ws.Cells.Item(99, 13).Value = -2748.4546  # M99: was -9261628
$ws.Cells.Item(99, 14).Value = -11296  # N99: was -9996
$ws.Cells.Item(126, 8).Value = 5513.1875  # H126: was 4907302
$ws.Cells.Item(126, 9).Value = 4246.4546  # I126: was 9263126
$ws.Cells.Item(126, 10).Value = 8300  # J126: was 7000
$ws.Cells.Item(126, 11).Value = 12739.3638  # K126: was 27789378
$ws.Cells.Item(126, 12).Value = 24900  # L126: was 21000
$ws.Cells.Item(126, 13).Value = -10269.3638  # M126: was -27786908
$ws.Cells.Item(126, 14).Value = -29840  # N126: was -25940

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 199  # H33: was 257.25
$ws.Cells.Item(33, 10).Value = 0  # J33: was 276.66666
$ws.Cells.Item(33, 12).Value = 0  # L33: was 1659.99996
$ws.Cells.Item(33, 14).ClearContents()  # N33: was -2225.99996
$ws.Cells.Item(68, 8).Value = 50501.5  # H68: was 100003
$ws.Cells.Item(68, 10).Value = 50501.5  # J68: was 100003
$ws.Cells.Item(68, 12).Value = 151504.5  # L68: was 300009
$ws.Cells.Item(68, 14).Value = -153126.5  # N68: was -301631
$ws.Cells.Item(71, 8).Value = 50501.5  # H71: was 100003
$ws.Cells.Item(71, 10).Value = 50501.5  # J71: was 100003
$ws.Cells.Item(71, 12).Value = 454513.5  # L71: was 900027
$ws.Cells.Item(71, 14).Value = -462625.5  # N71: was -908139
$ws.Cells.Item(92, 8).Value = 1098  # H92: was 898.5714
$ws.Cells.Item(92, 10).Value = 1945  # J92: was 1172.5
$ws.Cells.Item(92, 12).Value = 5835  # L92: was 3517.5
$ws.Cells.Item(92, 14).Value = -8331  # N92: was -6013.5
$ws.Cells.Item(97, 8).Value = 1081.25  # H97: was 807.1429000000001
$ws.Cells.Item(97, 10).Value = 1900  # J97: was 800
$ws.Cells.Item(97, 12).Value = 5700  # L97: was 2400
$ws.Cells.Item(97, 14).Value = -6692  # N97: was -3392
$ws.Cells.Item(131, 8).Value = 760.83  # H131: was 755.84
$ws.Cells.Item(131, 10).Value = 774.7732  # J131: was 769.6288500000001
$ws.Cells.Item(131, 12).Value = 2324.3196  # L131: was 2308.88655
$ws.Cells.Item(131, 14).Value = -12404.3196  # N131: was -12388.88655

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3941.2856  # H80: was 3321.1924
$ws.Cells.Item(80, 9).Value = 3496.6  # I80: was 2999.0833
$ws.Cells.Item(80, 10).Value = 5053  # J80: was 3597.2856
$ws.Cells.Item(80, 11).Value = 3496.6  # K80: was 2999.0833
$ws.Cells.Item(80, 12).Value = 5053  # L80: was 3597.2856
$ws.Cells.Item(80, 13).Value = -2498.6  # M80: was -2001.0833
$ws.Cells.Item(80, 14).Value = -7049  # N80: was -5593.2856
$ws.Cells.Item(83, 8).Value = 3941.2856  # H83: was 3321.1924
$ws.Cells.Item(83, 9).Value = 3496.6  # I83: was 2999.0833
$ws.Cells.Item(83, 10).Value = 5053  # J83: was 3597.2856
$ws.Cells.Item(83, 11).Value = 17483  # K83: was 14995.4165
$ws.Cells.Item(83, 12).Value = 25265  # L83: was 17986.428
$ws.Cells.Item(83, 13).Value = -12491  # M83: was -10003.4165
$ws.Cells.Item(83, 14).Value = -35249  # N83: was -27970.428
$ws.Cells.Item(107, 8).Value = 3183.8333  # H107: was 4150.75
$ws.Cells.Item(107, 10).Value = 3800.6  # J107: was 5501
$ws.Cells.Item(107, 12).Value = 3800.6  # L107: was 5501
$ws.Cells.Item(107, 14).Value = -7640.6  # N107: was -9341
$ws.Cells.Item(122, 8).Value = 1691.3636  # H122: was 1740.5
$ws.Cells.Item(122, 10).Value = 1667.6666  # J122: was 1761.2
$ws.Cells.Item(122, 12).Value = 5002.9998  # L122: was 5283.6
$ws.Cells.Item(122, 14).Value = -9902.9998  # N122: was -10183.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6588.8887  # H7: was 6108.6665
$ws.Cells.Item(7, 9).Value = 6216.6665  # I7: was 5700.4443
$ws.Cells.Item(7, 11).Value = 6216.6665  # K7: was 5700.4443
$ws.Cells.Item(7, 13).Value = -6104.6665  # M7: was -5588.4443
$ws.Cells.Item(22, 8).Value = 3750  # H22: was 2001.3846
$ws.Cells.Item(22, 9).Value = 3667  # I22: was 2775
$ws.Cells.Item(22, 10).Value = 3833  # J22: was 1657.5555
$ws.Cells.Item(22, 11).Value = 3667  # K22: was 2775
$ws.Cells.Item(22, 12).Value = 3833  # L22: was 1657.5555
$ws.Cells.Item(22, 13).Value = -3372  # M22: was -2480
$ws.Cells.Item(22, 14).Value = -4423  # N22: was -2247.5555
$ws.Cells.Item(27, 8).Value = 3750  # H27: was 2001.3846
$ws.Cells.Item(27, 9).Value = 3667  # I27: was 2775
$ws.Cells.Item(27, 10).Value = 3833  # J27: was 1657.5555
$ws.Cells.Item(27, 11).Value = 3667  # K27: was 2775
$ws.Cells.Item(27, 12).Value = 3833  # L27: was 1657.5555
$ws.Cells.Item(27, 13).Value = -3560  # M27: was -2668
$ws.Cells.Item(27, 14).Value = -4047  # N27: was -1871.5555
$ws.Cells.Item(46, 8).Value = 2186  # H46: was 2382.8572
$ws.Cells.Item(46, 9).Value = 1724.75  # I46: was 2000
$ws.Cells.Item(46, 10).Value = 2647.25  # J46: was 2670
$ws.Cells.Item(46, 11).Value = 1724.75  # K46: was 2000
$ws.Cells.Item(46, 12).Value = 2647.25  # L46: was 2670
$ws.Cells.Item(46, 13).Value = -1536.75  # M46: was -1812
$ws.Cells.Item(46, 14).Value = -3023.25  # N46: was -3046
$ws.Cells.Item(61, 8).Value = 4494.5  # H61: was 4684.048
$ws.Cells.Item(61, 9).Value = 1772.6154  # I61: was 1877.5
$ws.Cells.Item(61, 11).Value = 1772.6154  # K61: was 1877.5
$ws.Cells.Item(61, 13).Value = -1570.6154  # M61: was -1675.5
$ws.Cells.Item(82, 8).Value = 3255.5557  # H82: was 3640.5
$ws.Cells.Item(82, 9).Value = 4283.3335  # I82: was 3886
$ws.Cells.Item(82, 10).Value = 1200  # J82: was 3067.6667
$ws.Cells.Item(82, 11).Value = 4283.3335  # K82: was 3886
$ws.Cells.Item(82, 12).Value = 1200  # L82: was 3067.6667
$ws.Cells.Item(82, 13).Value = -3922.3335  # M82: was -3525
$ws.Cells.Item(82, 14).Value = -1922  # N82: was -3789.6667
$ws.Cells.Item(85, 8).Value = 3255.5557  # H85: was 3640.5
$ws.Cells.Item(85, 9).Value = 4283.3335  # I85: was 3886
$ws.Cells.Item(85, 10).Value = 1200  # J85: was 3067.6667
$ws.Cells.Item(85, 11).Value = 4283.3335  # K85: was 3886
$ws.Cells.Item(85, 12).Value = 1200  # L85: was 3067.6667
$ws.Cells.Item(85, 13).Value = -3035.3335  # M85: was -2638
$ws.Cells.Item(85, 14).Value = -3696  # N85: was -5563.6667
$ws.Cells.Item(113, 8).Value = 4494.5  # H113: was 4684.048
$ws.Cells.Item(113, 9).Value = 1772.6154  # I113: was 1877.5
$ws.Cells.Item(113, 11).Value = 1772.6154  # K113: was 1877.5
$ws.Cells.Item(113, 13).Value = 397.3846000000001  # M113: was 292.5
$ws.Cells.Item(126, 8).Value = 6588.8887  # H126: was 6108.6665
$ws.Cells.Item(126, 9).Value = 6216.6665  # I126: was 5700.4443
$ws.Cells.Item(126, 11).Value = 18649.9995  # K126: was 17101.3329
$ws.Cells.Item(126, 13).Value = -16179.9995  # M126: was -14631.3329
$ws.Cells.Item(140, 8).Value = 48694.25  # H140: was 49241.668
$ws.Cells.Item(140, 10).Value = 48694.25  # J140: was 49241.668
$ws.Cells.Item(140, 12).Value = 48694.25  # L140: was 49241.668
$ws.Cells.Item(140, 14).Value = -59054.25  # N140: was -59601.668

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 3247960.5  # H107: was 3497788.2
$ws.Cells.Item(107, 9).Value = 399.83334  # I107: was 440
$ws.Cells.Item(107, 11).Value = 1199.50002  # K107: was 1320
$ws.Cells.Item(107, 13).Value = 720.4999800000001  # M107: was 600
$ws.Cells.Item(108, 8).Value = 0  # H108: was 32000
$ws.Cells.Item(108, 10).Value = 0  # J108: was 32000
$ws.Cells.Item(108, 12).Value = 0  # L108: was 32000
$ws.Cells.Item(108, 14).ClearContents()  # N108: was -39680
$ws.Cells.Item(126, 8).Value = 1929.75  # H126: was 1932.25
$ws.Cells.Item(126, 9).Value = 1047.6  # I126: was 1051.6
$ws.Cells.Item(126, 11).Value = 3142.8  # K126: was 3154.8
$ws.Cells.Item(126, 13).Value = -672.7999999999997  # M126: was -684.7999999999997
